$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Formula = "'123"

$ws.Range("C2").Value2 = [double]"8.411250374824416e-10"
$ws.Range("D2").Value2 = [double]"4.527853396086863e-08"
$ws.Range("E2").Value2 = [double]"7.285874819458044e-08"
$ws.Range("F2").Value2 = [double]"1.510857908510606e-23"
$ws.Range("G2").Value2 = [double]"4.153595805285083e-16"
$ws.Range("H2").Value2 = [double]"3.086748091649903e-09"
$ws.Range("I2").Value2 = [double]"0.3692929004708708"
$ws.Range("J2").Value2 = [double]"0.0002954943620861417"
$ws.Range("K2").Value2 = [double]"0.005014276079098002"
$ws.Range("M2").Value2 = [double]"0.01989595693847302"
$ws.Range("N2").Value2 = [double]"99.57449200792134"
$ws.Range("O2").Value2 = [double]"0.01656385322649616"
$ws.Range("P2").Value2 = [double]"0.003915020881706073"
$ws.Range("Q2").Value2 = [double]"0.0006615052573192528"
$ws.Range("R2").Value2 = [double]"0.003663183279454354"
$ws.Range("S2").Value2 = [double]"0.002920768697589791"
$ws.Range("T2").Value2 = [double]"0.003275591990947319"
$ws.Range("U2").Value2 = [double]"8.359851092386573e-06"
$ws.Range("V2").Value2 = [double]"9.443867601533344e-07"
$ws.Range("W2").Value2 = [double]"4.547602998846134e-12"
$ws.Range("X2").Value2 = [double]"1.830181898234014e-11"
$ws.Range("Y2").Value2 = [double]"1.043903303088885e-14"
$ws.Range("Z2").Value2 = [double]"1.125559005260744e-11"
$ws.Range("AA2").Value2 = [double]"1.718676837236395e-12"
$ws.Range("AB2").Value2 = [double]"4.064992702953596e-14"
$ws.Range("AC2").Value2 = [double]"4.147929409899162e-13"
$ws.Range("AD2").Value2 = [double]"8.589114370037159e-16"
$ws.Range("AE2").Value2 = [double]"5.286966308310837e-17"
$ws.Range("AF2").Value2 = [double]"2.296873135697792e-19"
$ws.Range("AG2").Value2 = [double]"1.117756785894828e-20"
$ws.Range("AH2").Value2 = [double]"3.095668286238961e-21"
$ws.Range("AI2").Value2 = [double]"4.731452080357051e-22"
$ws.Range("AJ2").Value2 = [double]"2.693750671942818e-22"
$ws.Range("AL2").Value2 = [double]"1.26689659687138e-08"
$ws.Range("AM2").Value2 = [double]"1.557938990266146e-09"
$ws.Range("AN2").Value2 = [double]"3.284285956940382e-10"
$ws.Range("AR2").Value2 = [double]"1.100103227258436e-18"
$ws.Range("AS2").Value2 = [double]"1.387335337860952e-16"
$ws.Range("AT2").Value2 = [double]"1.986282012568797e-16"
$ws.Range("AU2").Value2 = [double]"5.565307756481919e-32"
$ws.Range("AV2").Value2 = [double]"9.985072070248802e-25"
$ws.Range("AW2").Value2 = [double]"9.669133692984975e-18"
$ws.Range("AX2").Value2 = [double]"1.11866273983658e-09"
$ws.Range("AY2").Value2 = [double]"2.960124172061386e-12"
$ws.Range("AZ2").Value2 = [double]"5.360772877655186e-10"
$ws.Range("BB2").Value2 = [double]"2.657854801722107e-07"
$ws.Range("BC2").Value2 = [double]"0.01990732032944903"
$ws.Range("BD2").Value2 = [double]"3.293278795176256e-05"
$ws.Range("BE2").Value2 = [double]"0.01647298272391547"
$ws.Range("BF2").Value2 = [double]"0.01954603259840981"
$ws.Range("BG2").Value2 = [double]"2.745467814459453"
$ws.Range("BH2").Value2 = [double]"16.62491709878523"
$ws.Range("BI2").Value2 = [double]"36.1250022581795"
$ws.Range("BJ2").Value2 = [double]"14.63680922797722"
$ws.Range("BK2").Value2 = [double]"0.03138631167399627"
$ws.Range("BL2").Value2 = [double]"0.01232703299549507"
$ws.Range("BM2").Value2 = [double]"0.01761354375923693"
$ws.Range("BN2").Value2 = [double]"0.01438153742321458"
$ws.Range("BO2").Value2 = [double]"2.015469985737537"
$ws.Range("BP2").Value2 = [double]"1.056015778528811"
$ws.Range("BQ2").Value2 = [double]"5.349932255897751"
$ws.Range("BR2").Value2 = [double]"0.5937520558127614"
$ws.Range("BS2").Value2 = [double]"13.6645142006224"
$ws.Range("BT2").Value2 = [double]"0.001033681983592556"
$ws.Range("BU2").Value2 = [double]"2.181884750703781"
$ws.Range("BV2").Value2 = [double]"0.07190768514262993"
$ws.Range("BW2").Value2 = [double]"0.003888433489884459"
$ws.Range("BX2").Value2 = [double]"0.04772044302472385"
$ws.Range("BY2").Value2 = [double]"1.847000158650691"
$ws.Range("CA2").Value2 = [double]"2.124358706248295"
$ws.Range("CB2").Value2 = [double]"0.5444438823336096"
$ws.Range("CC2").Value2 = [double]"0.2342136206872819"
